# Applies the cryptos.xlsx price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.349.54'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = '1.833.60'
$ws.Range('E3').Value = '  +1.05%  '
$ws.Range('E4').Value = '  +0.91%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.91'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.86%  '
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4750'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3688'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07457'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8862'
$ws.Range('D10').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.45'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.07%  '
$ws.Range('D12').Value = '1.870.45'
$ws.Range('E12').Value = '  +1.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07324'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +3.06%  '
$ws.Range('E14').Value = '  +1.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.20'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.28%  '
$ws.Range('E16').Value = '  +1.17%  '
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008805'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.011'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.80'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.23%  '
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').Value = '27.504.65'
$ws.Range('E21').Value = '  +2.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.292'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  +0.91%  '
$ws.Range('D24').Value = '2.089.30'
$ws.Range('E24').Value = '  +1.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.897'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.95'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('E27').Value = '  +1.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.143'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.243'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '117.38'
$ws.Range('D30').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08995'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7536'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('E33').Value = '  +1.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.548'
$ws.Range('D34').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.944'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.44%  '
$ws.Range('E36').Value = '  +0.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.104'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05342'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01956'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.64%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.979'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.276'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.395'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5317'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.59%  '
$ws.Range('E44').Value = '  +0.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.482'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.92%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4913'
$ws.Range('D46').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.55'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.03%  '
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('E49').Value = '  +0.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.672'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06296'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.08%  '
